# Excel COM-interop script: add a "Player Info" sheet, and on both "ODI
# Batting" and "ODI Bowling" rename the MATCH_CARD_LINK column to
# MATCH_CODE, collapsing each scorecard URL down to just the bare numeric
# match code.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Player Info" sheet in front of everything else ---
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Header row: bold, thin box border, centered/top-aligned - matches the
# styling already used for header rows on the other sheets.
$hdr = $playerInfo.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Data row. The ID is text (matches the rest of the workbook, where every
# "numeric" value is stored as text) - use a leading apostrophe so Excel
# keeps it as text instead of inferring a number, then reset the style so
# the quote-prefix flag doesn't leave a stray format behind.
$playerInfo.Range("A2").Value = "'4383"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Andrew James Tye"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# --- 2. "ODI Batting": MATCH_CARD_LINK (col D) -> MATCH_CODE, URL -> code ---
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @("4108", "4115", "4123", "4125", "4166", "4167", "4168")
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    Set-TextValue $batting.Range("D$row") $battingCodes[$i]
}

# --- 3. "ODI Bowling": MATCH_CARD_LINK (col B) -> MATCH_CODE, URL -> code ---
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @("4108", "4115", "4123", "4125", "4166", "4167", "4168")
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $row = $i + 2
    Set-TextValue $bowling.Range("B$row") $bowlingCodes[$i]
}
